$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.238.39"
$ws.Range("E2").Value = "  +0.00%  "

$ws.Range("D3").Value = "1.603.73"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").Value = "'212.22"

$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("E8").Value = "  -0.29%  "

$ws.Range("E9").Value = "  -0.69%  "

$ws.Range("D10").Value = "'18.17"
$ws.Range("E10").Value = "  +0.19%  "

$ws.Range("D11").Value = "'0.0814"
$ws.Range("E11").Value = "  -0.36%  "

$ws.Range("D12").Value = "1.824.67"
$ws.Range("E12").Value = "  -0.50%  "

$ws.Range("D13").Value = "1.607.45"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("E15").Value = "  +0.75%  "

$ws.Range("D16").Value = "26.219.58"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "'61.38"
$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("D20").Value = "'202.45"
$ws.Range("E20").Value = "  +1.76%  "

$ws.Range("E21").Value = "  +1.01%  "

$ws.Range("E22").Value = "  -1.58%  "

$ws.Range("D23").Value = "'5.98"
$ws.Range("E23").Value = "  -0.62%  "

$ws.Range("E24").Value = "  +9.61%  "

$ws.Range("D25").Value = "'144.01"
$ws.Range("E25").Value = "  +0.84%  "

$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("D27").Value = "'0.121"
$ws.Range("E27").Value = "  -6.79%  "

$ws.Range("D28").Value = "'15.20"
$ws.Range("E28").Value = "  +0.18%  "

$ws.Range("D29").Value = "'6.56"
$ws.Range("E29").Value = "  +1.12%  "

$ws.Range("E30").Value = "  +3.24%  "

$ws.Range("E31").Value = "  -0.41%  "

$ws.Range("E32").Value = "  +2.08%  "

$ws.Range("E33").Value = "  -2.66%  "

$ws.Range("E34").Value = "  +2.89%  "

$ws.Range("E35").Value = "  -0.51%  "

$ws.Range("D36").Value = "1.146.68"
$ws.Range("E36").Value = "  +3.41%  "

$ws.Range("E37").Value = "  +7.97%  "

$ws.Range("E38").Value = "  -0.12%  "

$ws.Range("D39").Value = "'0.795"
$ws.Range("E39").Value = "  +1.23%  "

$ws.Range("E40").Value = "  -0.17%  "

$ws.Range("D41").Value = "'0.498"
$ws.Range("E41").Value = "  -0.64%  "

$ws.Range("D42").Value = "'0.783"
$ws.Range("E42").Value = "  +0.36%  "

$ws.Range("E43").Value = "  +2.70%  "

$ws.Range("D44").Value = "1.738.50"
$ws.Range("E44").Value = "  -0.41%  "

$ws.Range("D45").Value = "'91.70"
$ws.Range("E45").Value = "  -0.87%  "

$ws.Range("E46").Value = "  -2.01%  "

$ws.Range("D47").Value = "'54.03"
$ws.Range("E47").Value = "  +0.49%  "

$ws.Range("D48").Value = "'0.0505"
$ws.Range("E48").Value = "  -0.66%  "

# Row 49-51 coin rotation: BabyDogeCoin moves to 49, Mantle moves to 50, USDD moves to 51
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₇0968"
$ws.Range("E49").Value = "  -10.22%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.407"
$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.17%  "
